# edit.ps1
# Applies two changes to RBTV.docx:
#  1. Removes the empty paragraph that immediately follows the "RBTV" title
#     paragraph at the top of the document.
#  2. Adds a new bulleted list item with the text
#     "Khách hang chỉ được đặt món ăn trong thực đơn" right after the last
#     existing populated bulleted item (the one ending in "... thuế chân."),
#     i.e. right before the trailing empty bullet paragraph at the end of
#     the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: delete the empty paragraph right after the "RBTV" paragraph.
# ---------------------------------------------------------------------
$titleParaIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "RBTV") {
        $titleParaIndex = $i
        break
    }
}

$emptyParaIndex = $titleParaIndex + 1
$emptyPara = $d.Paragraphs.Item($emptyParaIndex)
if ($emptyPara.Range.Text.Trim() -eq "") {
    $emptyPara.Range.Delete()
}

# ---------------------------------------------------------------------
# Change 2: insert a new bullet paragraph after the "... thuế chân."
# paragraph (the last populated ListParagraph bullet item), before the
# final, already-empty bullet paragraph at the end of the document.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("thuế chân", $true, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)

if ($found) {
    $matchEnd = $rng.End
    $anchorParaIndex = $d.Range(0, $matchEnd).Paragraphs.Count
} else {
    # Fallback: anchor on the paragraph just before the final (empty) one.
    $anchorParaIndex = $d.Paragraphs.Count - 1
}

$anchorPara = $d.Paragraphs.Item($anchorParaIndex)
$anchorPara.Range.InsertParagraphAfter()

$newParaIndex = $anchorParaIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Range.Text = "Khách hang chỉ được đặt món ăn trong thực đơn"
